# Fruta / hortaliza, semanal
# New weekly price-report rows are inserted for Chirimoya (La Palmera de La Serena).
# A new set of 3 rows (Especial/Primera/Segunda, Provincia de Limarí) for the week
# of 2021-11-18 is inserted before the existing row 81, pushing the rest of the
# table (rows 81-113) down to rows 84-116.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at position 81 (shifts old rows 81:113 down to 84:116)
$ws.Rows("81:83").Insert()

# Constant values shared by every data row in this sheet
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100107
$producto  = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad  = "Cultivar IV Región"
$unidad    = "$/kilo (en caja de 15 kilos)"
$origen    = "Provincia de Limarí"

# New row 81: Especial
$r = 81
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value   = $mercado
$ws.Cells.Item($r, 3).Value   = $region
$ws.Cells.Item($r, 4).Value2  = 44518
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value   = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value   = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value  = $categoria
$ws.Cells.Item($r, 11).Value  = $variedad
$ws.Cells.Item($r, 12).Value  = "Especial"
$ws.Cells.Item($r, 13).Value2 = 400
$ws.Cells.Item($r, 14).Value2 = 1800
$ws.Cells.Item($r, 15).Value2 = 1900
$ws.Cells.Item($r, 16).Value2 = 1850
$ws.Cells.Item($r, 17).Value  = $unidad
$ws.Cells.Item($r, 18).Value  = $origen
$ws.Cells.Item($r, 19).Value2 = 1850
$ws.Cells.Item($r, 20).Value2 = 1

# New row 82: Primera
$r = 82
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value   = $mercado
$ws.Cells.Item($r, 3).Value   = $region
$ws.Cells.Item($r, 4).Value2  = 44518
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value   = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value   = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value  = $categoria
$ws.Cells.Item($r, 11).Value  = $variedad
$ws.Cells.Item($r, 12).Value  = "Primera"
$ws.Cells.Item($r, 13).Value2 = 400
$ws.Cells.Item($r, 14).Value2 = 1500
$ws.Cells.Item($r, 15).Value2 = 1600
$ws.Cells.Item($r, 16).Value2 = 1550
$ws.Cells.Item($r, 17).Value  = $unidad
$ws.Cells.Item($r, 18).Value  = $origen
$ws.Cells.Item($r, 19).Value2 = 1550
$ws.Cells.Item($r, 20).Value2 = 1

# New row 83: Segunda
$r = 83
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value   = $mercado
$ws.Cells.Item($r, 3).Value   = $region
$ws.Cells.Item($r, 4).Value2  = 44518
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value   = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value   = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value  = $categoria
$ws.Cells.Item($r, 11).Value  = $variedad
$ws.Cells.Item($r, 12).Value  = "Segunda"
$ws.Cells.Item($r, 13).Value2 = 300
$ws.Cells.Item($r, 14).Value2 = 1200
$ws.Cells.Item($r, 15).Value2 = 1300
$ws.Cells.Item($r, 16).Value2 = 1250
$ws.Cells.Item($r, 17).Value  = $unidad
$ws.Cells.Item($r, 18).Value  = $origen
$ws.Cells.Item($r, 19).Value2 = 1250
$ws.Cells.Item($r, 20).Value2 = 1
